$wb = $excel.ActiveWorkbook

# --- Sheet "Year of Examination" (sheet1) ---
# Remove the "Visual Variable" value in B4 (was MAIN-COL-DISC) - attribute no
# longer assigned to this sheet.
$wsYear = $wb.Worksheets.Item("Year of Examination")
$wsYear.Range("B4").ClearContents()

# --- Sheet "Disease Free Survival" (sheet2) ---
# Re-assign the Visual Variable from MAIN-COL-DISC to HEIGHT.
$wsDFS = $wb.Worksheets.Item("Disease Free Survival")
$wsDFS.Range("B4").Value = "HEIGHT"

# --- Sheet "Staging T" (sheet3) ---
# Re-assign the Visual Variable from BASE-COL-DISC to MAIN-COL-DISC, and add
# the new discretization-weight columns E/F/G for rows 2-7.
$wsT = $wb.Worksheets.Item("Staging T")
$wsT.Range("B4").Value = "MAIN-COL-DISC"

$wsT.Range("E2").Value = 0.15
$wsT.Range("F2").Value = 0.15
$wsT.Range("G2").Value = 0.15

$wsT.Range("E3").Value = 0.15
$wsT.Range("F3").Value = 0.15
$wsT.Range("G3").Value = 0.15

$wsT.Range("E4").Value = 0.3
$wsT.Range("F4").Value = 0.3
$wsT.Range("G4").Value = 0.3

$wsT.Range("E5").Value = 0
$wsT.Range("F5").Value = 0.3
$wsT.Range("G5").Value = 0.3

$wsT.Range("E6").Value = 0.3
$wsT.Range("F6").Value = 0.1
$wsT.Range("G6").Value = 0.1

$wsT.Range("E7").Value = 0.6
$wsT.Range("F7").Value = 0
$wsT.Range("G7").Value = 0

# --- Sheet "Staging N" (sheet4) ---
# Re-assign the Visual Variable from MAIN-COL-DISC to BASE-COL-DISC.
$wsN = $wb.Worksheets.Item("Staging N")
$wsN.Range("B4").Value = "BASE-COL-DISC"

# --- Sheet "Sex" (sheet5) ---
# Flip Visual Mapping from Yes to No, and clear the now-unused Visual
# Variable cell.
$wsSex = $wb.Worksheets.Item("Sex")
$wsSex.Range("B3").Value = "No"
$wsSex.Range("B4").ClearContents()

# --- Selections / active sheet ---
# Set each sheet's remembered selection, finishing with "Staging T" as the
# active tab (matches the saved workbook view state).
$wsYear.Activate()
$wsYear.Range("B5").Select()

$wsDFS.Activate()
$wsDFS.Range("B4").Select()

$wsT.Activate()
$wsT.Range("B4").Select()

$wsN.Activate()
$wsN.Range("B4").Select()

$wsSex.Activate()
$wsSex.Range("B14").Select()

$wsT.Activate()
